$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 61: change border styling from the old "last row" xf slot to the
# shared bottom-border slot already used by row 59 (C:H). Visually identical
# (still a thin bottom border); only the underlying style-table slot moves so
# that the slot it vacates can be reused (without a border) by the new last
# row below it.
$row61 = $ws.Range("C61:H61")
$row61.NumberFormat = $ws.Range("C59").NumberFormat
$row61.Borders.Item(9).LineStyle = 1
$ws.Range("C61").NumberFormat = "00"
$ws.Range("E61").NumberFormat = "00"
$ws.Range("G61").NumberFormat = "00"

# --- Row 62: new chapter entry ("Working with tables" / "Using CALCULATETABLE").
$ws.Range("C62").Value = 12
$ws.Range("C62").NumberFormat = "00"

$ws.Range("D62").Value = "Working with tables"
$ws.Range("D62").NumberFormat = "General"

$ws.Range("E62").Value = 1
$ws.Range("E62").NumberFormat = "00"

$ws.Range("F62").Value = "Using CALCULATETABLE"
$ws.Range("F62").NumberFormat = "General"

$ws.Range("G62").Value = 1
$ws.Range("G62").NumberFormat = "00"

$ws.Range("H62").Value = "Using CALCULATETABLE"
$ws.Range("H62").NumberFormat = "General"

$ws.Range("I62").Value = "CALCULATETABLE"
$ws.Range("I62").NumberFormat = "00"

$ws.Range("B62").Formula = "=CONCAT(TEXT(C62,""00""),TEXT(E62,""00""),TEXT(G62,""00""))"

# New row has no border (it's now the open-ended last row of the table).
$ws.Range("C62:I62").Borders.Item(9).LineStyle = -4142

# Move the selection the way the author's session ended up.
$ws.Range("I62").Select()
